$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$words = @(
    'kümmern',
    'enden',
    'scheinen',
    'wirken',
    'formen',
    'heben',
    'äußern',
    'zielen',
    'hören',
    'bitten',
    'betteln',
    'fällen',
    'lesen',
    'pfeifen',
    'schulden',
    'werden',
    'werfen',
    'flüchten',
    'wehtun',
    'knarren',
    'saufen',
    'kehren',
    'runden',
    'dringen',
    'geben',
    'feiern',
    'stammen',
    'boxen',
    'biegen',
    'fangen',
    'sperren',
    'schrecken',
    'sinken',
    'lügen',
    'zögern',
    'schwächen',
    'schreiten',
    'stehlen',
    'brauchen',
    'graben',
    'zünden',
    'fließen',
    'kichern',
    'bellen',
    'klingen',
    'planen',
    'schmecken',
    'greifen',
    'liefern',
    'fallen',
    'heilen',
    'tollen',
    'decken',
    'münzen',
    'rasen',
    'scheitern',
    'mauern',
    'malen',
    'weichen',
    'streichen',
    'machen',
    'suchen',
    'achten',
    'altern',
    'bluten',
    'jagen',
    'heulen',
    'ändern',
    'wachsen',
    'schlucken',
    'wundern',
    'pflanzen',
    'flehen',
    'spielen',
    'spüren',
    'ehren',
    'dienen',
    'seufzen',
    'gelten',
    'schwingen',
    'führen',
    'sterben',
    'filmen',
    'treiben',
    'liegen',
    'sprengen',
    'mögen',
    'lockern',
    'schenken',
    'töten',
    'gründen',
    'jubeln',
    'irren',
    'fahren',
    'hauen',
    'wüten',
    'warnen',
    'reizen',
    'trauen',
    'helfen',
    'backen',
    'loben',
    'wenden',
    'platzen',
    'räumen',
    'arten',
    'klettern',
    'stecken',
    'schwören',
    'ärgern',
    'quälen',
    'siegen',
    'rufen',
    'bergen',
    'kosten',
    'trennen',
    'bauen',
    'fischen',
    'tropfen',
    'freuen',
    'erben',
    'grüßen',
    'drehen',
    'sichern',
    'sorgen',
    'zeigen',
    'folgen',
    'spinnen',
)

for ($i = 0; $i -lt $words.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $words[$i]
}

